$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy H1's formatting (bold, bordered, centered header style) to I1/J1, then set text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Data rows 2-15: I = 1, J = same value as H
for ($r = 2; $r -le 15; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $h
}
